$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Translate German publisher/contact text to English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Fill in the previously empty Description value cell
$ws.Range("B12").Value = "Possible match statuses in the context of adding patient identities. "
